$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - cán bộ "Trần Quốc T"
$ws.Range("A2").Value = "D:\images\user.png"
$ws.Range("B2").Value = "cb43"
$ws.Range("C2").Value = "Trần Quốc T"
$ws.Range("D2").Value = "Lĩnh vực tài nguyên và đất đai"

# Row 3 - cán bộ "Hà Anh T"
$ws.Range("A3").Value = "D:\images\user.png"
$ws.Range("B3").Value = "cb44"
$ws.Range("C3").Value = "Hà Anh T"
$ws.Range("D3").Value = "Lĩnh vực xây dựng"

$ws.PageSetup.Orientation = 1

$ws.Range("A3").Select() | Out-Null
